$d = $word.ActiveDocument
$matches = 0

# Phase 1: replace each original run's full text with a unique placeholder token
# to avoid any collisions between old/new text blocks that overlap across paragraphs.
$ok = $d.Content.Find.Execute("Introduzir o aluno na engenharia das reações químicas, através dos conceitos fundamentais da cinética química aplicada a reatores químicos ideais.", $true, $false, $false, $false, $false, $true, 1, $false, "@@BLOCK0@@", 2)
if (-not $ok) { Write-Output "PHASE1 FAILED at index 0 (para=5 run=0)" } else { $matches = $matches + 1 }
$ok = $d.Content.Find.Execute("Introduction to Chemical Reaction Engineering through the fundamental concepts of chemical kinetics applied to ideal chemical reactors", $true, $false, $false, $false, $false, $true, 1, $false, "@@BLOCK1@@", 2)
if (-not $ok) { Write-Output "PHASE1 FAILED at index 1 (para=6 run=0)" } else { $matches = $matches + 1 }
$ok = $d.Content.Find.Execute("5963230 - Leandro Gonçalves de Aguiar^l", $true, $false, $false, $false, $false, $true, 1, $false, "@@BLOCK2@@", 2)
if (-not $ok) { Write-Output "PHASE1 FAILED at index 2 (para=8 run=0)" } else { $matches = $matches + 1 }
$ok = $d.Content.Find.Execute("6310316 - Liana Alvares Rodrigues", $true, $false, $false, $false, $false, $true, 1, $false, "@@BLOCK3@@", 2)
if (-not $ok) { Write-Output "PHASE1 FAILED at index 3 (para=8 run=1)" } else { $matches = $matches + 1 }
$ok = $d.Content.Find.Execute("1. Introdução a cinética. 2. Reações a volume constante. 3. Reações a volume variável. 4. Modelos ideais de reatores químicos isotérmicos. 5. Análise de dados cinéticos em reatores químicos isotérmicos", $true, $false, $false, $false, $false, $true, 1, $false, "@@BLOCK4@@", 2)
if (-not $ok) { Write-Output "PHASE1 FAILED at index 4 (para=10 run=0)" } else { $matches = $matches + 1 }
$ok = $d.Content.Find.Execute("1. Introduction to Kinetics. 2. Reactions at Constant Volume. 3. Reactions at Variable Volume. 4. Ideal Models of Isothermal Chemical Reactors. 5. Analysis of Kinetic Data in Isothermal Chemical Reactors", $true, $false, $false, $false, $false, $true, 1, $false, "@@BLOCK5@@", 2)
if (-not $ok) { Write-Output "PHASE1 FAILED at index 5 (para=11 run=0)" } else { $matches = $matches + 1 }
$ok = $d.Content.Find.Execute("1. INTRODUÇÃO A CINÉTICA^lTipos de Reações Químicas. Lei de velocidade e seus principais parâmetros. Influência da temperatura sobre a taxa da reação. Ativação das reações químicas Equação de Arrhenius. Energia de ativação. Conversão. Concentração e sua variação numa transformação química. ^l2. REAÇÕES A VOLUME CONSTANTE^lReações irreversíveis de ordem um. Reações irreversíveis de ordem dois. Reações irreversíveis de ordem três. Reações irreversíveis de ordem qualquer. ^l3. REAÇÕES A VOLUME VARIÁVEL^lConceitos. Fração de conversão volumétrica. Reações a volume variável de ordem um e dois. ^l4. MODELOS IDEAIS DE REATORES QUÍMICOS ISOTÉRMICOS: ^lEquações fundamentais de projeto de reatores. Reator tanque descontínuo (BSTR). Reator tanque de mistura contínuo (CSTR). Reator tubular de fluxo pistonado (PFR). Comparação de desempenho de reatores CSTR e PFR. Reatores CSTR em cascata. Associação mista de reatores em série: CSTR e PFR ^l5. ANÁLISE DE DADOS CINÉTICOS EM REATORES QUÍMICOS ISOTÉRMICOS^lBalanço de massa e coleta de dados em reatores ideais isotérmicos: batelada (BSTR), reator tanque de mistura contínuo (CSTR) e Reator tubular (PFR)", $true, $false, $false, $false, $false, $true, 1, $false, "@@BLOCK6@@", 2)
if (-not $ok) { Write-Output "PHASE1 FAILED at index 6 (para=13 run=0)" } else { $matches = $matches + 1 }
$ok = $d.Content.Find.Execute("Duas provas escritas (P1 e P2) e eventuais trabalhos relacionados à disciplina^l", $true, $false, $false, $false, $false, $true, 1, $false, "@@BLOCK7@@", 2)
if (-not $ok) { Write-Output "PHASE1 FAILED at index 7 (para=16 run=1)" } else { $matches = $matches + 1 }
$ok = $d.Content.Find.Execute("Média da Primeira Avaliação (N) = 50% P1 + 50% P2.^lObs: fica a critério de cada docente a inserção de trabalhos no decorrer do curso, bem como a alteração do peso de cada prova em decorrência dos mesmos.^l", $true, $false, $false, $false, $false, $true, 1, $false, "@@BLOCK8@@", 2)
if (-not $ok) { Write-Output "PHASE1 FAILED at index 8 (para=16 run=3)" } else { $matches = $matches + 1 }
$ok = $d.Content.Find.Execute("Média Final = (N + Prova Recuperação)/2", $true, $false, $false, $false, $false, $true, 1, $false, "@@BLOCK9@@", 2)
if (-not $ok) { Write-Output "PHASE1 FAILED at index 9 (para=16 run=5)" } else { $matches = $matches + 1 }
$ok = $d.Content.Find.Execute("1- FOGLER, H.S. Elementos de engenharia das reações químicas. 3.ed. Rio de Janeiro: LTC Editora, 2009.^l^l2- LEVENSPIEL, O. Engenharia Das Reações Químicas, E ed (Blucher, São Paulo, 2000)^l3- VAN SANTEN, R.A.; Niemantsverdriet, J.W. Chemical kinetics and catalysis. New York: Plenum Press, 1995.^l4- Missen, R.W.; Mims, C.A.; Saville, B.A. Introduction to chemical reaction engineering and kinetics. New York: J. Wiley, 1999.^l5- Rothenberg, G. Catalysis: concepts and green applications. Weinheim: Wiley-VCH, 2008 Chichester.^l6- DENISOV, E.T.; Sarkisov, O.M.; Likhtenshtein, G.I. Chemical kinetics: fundamentals and new developments. Amsterdam: Elsevier, 2003.^l7- Hagen, J. Industrial catalysis: a practical approach. Weinheim: Wiley-VCH, 2006.^l8- Salmi, T.O.; Mikkola, J.; Warna, J.P. Chemical reaction engineering and reactor technology. Boca Raton: CRC Press/Taylor & Francis, 2011.^l9- Mortimer, M.; Taylor, P.G. Chemical kinetics and mechanism. Cambridge: Royal Society of Chemistry, 2002.^l10- FROMENT, G.F.; BISCHOFF, K.B. Chemical reactor analysis and design. 2nd. Ed. New York: John Wiley & Sons, 1990.^l11- HILL, C.G. An Introduction to chemical engineering kinetics and reactor design. New York: John Wiley&Sons, 1977.^l12- SMITH, J.M. Chemical engineering kinetics. 3rd. ed New York: McGraw-Hill,1981.^l13- DENBIGH, K.; TURNER, R. Introduction to chemical Reaction Design. Cambridge: Cambridge University Press, 1970.^l14 - AGUIAR, L. G. Problemas de cinética e reatores químicos. Curitiba: Appris Editora, 2023.", $true, $false, $false, $false, $false, $true, 1, $false, "@@BLOCK10@@", 2)
if (-not $ok) { Write-Output "PHASE1 FAILED at index 10 (para=18 run=0)" } else { $matches = $matches + 1 }

# Phase 2: replace each placeholder token with its final text
$ok = $d.Content.Find.Execute("@@BLOCK0@@", $true, $false, $false, $false, $false, $true, 1, $false, "1. Introdução a cinética. 2. Reações a volume constante. 3. Reações a volume variável. 4. Modelos ideais de reatores químicos isotérmicos. 5. Análise de dados cinéticos em reatores químicos isotérmicos", 2)
if (-not $ok) { Write-Output "PHASE2 FAILED at index 0 (para=5 run=0)" }
$ok = $d.Content.Find.Execute("@@BLOCK1@@", $true, $false, $false, $false, $false, $true, 1, $false, "1. Introduction to Kinetics. 2. Reactions at Constant Volume. 3. Reactions at Variable Volume. 4. Ideal Models of Isothermal Chemical Reactors. 5. Analysis of Kinetic Data in Isothermal Chemical Reactors", 2)
if (-not $ok) { Write-Output "PHASE2 FAILED at index 1 (para=6 run=0)" }
$ok = $d.Content.Find.Execute("@@BLOCK2@@", $true, $false, $false, $false, $false, $true, 1, $false, "Introduzir o aluno na engenharia das reações químicas, através dos conceitos fundamentais da cinética química aplicada a reatores químicos ideais.^l", 2)
if (-not $ok) { Write-Output "PHASE2 FAILED at index 2 (para=8 run=0)" }
$ok = $d.Content.Find.Execute("@@BLOCK3@@", $true, $false, $false, $false, $false, $true, 1, $false, "1. INTRODUÇÃO A CINÉTICA^lTipos de Reações Químicas. Lei de velocidade e seus principais parâmetros. Influência da temperatura sobre a taxa da reação. Ativação das reações químicas Equação de Arrhenius. Energia de ativação. Conversão. Concentração e sua variação numa transformação química. ^l2. REAÇÕES A VOLUME CONSTANTE^lReações irreversíveis de ordem um. Reações irreversíveis de ordem dois. Reações irreversíveis de ordem três. Reações irreversíveis de ordem qualquer. ^l3. REAÇÕES A VOLUME VARIÁVEL^lConceitos. Fração de conversão volumétrica. Reações a volume variável de ordem um e dois. ^l4. MODELOS IDEAIS DE REATORES QUÍMICOS ISOTÉRMICOS: ^lEquações fundamentais de projeto de reatores. Reator tanque descontínuo (BSTR). Reator tanque de mistura contínuo (CSTR). Reator tubular de fluxo pistonado (PFR). Comparação de desempenho de reatores CSTR e PFR. Reatores CSTR em cascata. Associação mista de reatores em série: CSTR e PFR ^l5. ANÁLISE DE DADOS CINÉTICOS EM REATORES QUÍMICOS ISOTÉRMICOS^lBalanço de massa e coleta de dados em reatores ideais isotérmicos: batelada (BSTR), reator tanque de mistura contínuo (CSTR) e Reator tubular (PFR)", 2)
if (-not $ok) { Write-Output "PHASE2 FAILED at index 3 (para=8 run=1)" }
$ok = $d.Content.Find.Execute("@@BLOCK4@@", $true, $false, $false, $false, $false, $true, 1, $false, "Duas provas escritas (P1 e P2) e eventuais trabalhos relacionados à disciplina", 2)
if (-not $ok) { Write-Output "PHASE2 FAILED at index 4 (para=10 run=0)" }
$ok = $d.Content.Find.Execute("@@BLOCK5@@", $true, $false, $false, $false, $false, $true, 1, $false, "Introduction to Chemical Reaction Engineering through the fundamental concepts of chemical kinetics applied to ideal chemical reactors", 2)
if (-not $ok) { Write-Output "PHASE2 FAILED at index 5 (para=11 run=0)" }
$ok = $d.Content.Find.Execute("@@BLOCK6@@", $true, $false, $false, $false, $false, $true, 1, $false, "Média da Primeira Avaliação (N) = 50% P1 + 50% P2.^lObs: fica a critério de cada docente a inserção de trabalhos no decorrer do curso, bem como a alteração do peso de cada prova em decorrência dos mesmos.", 2)
if (-not $ok) { Write-Output "PHASE2 FAILED at index 6 (para=13 run=0)" }
$ok = $d.Content.Find.Execute("@@BLOCK7@@", $true, $false, $false, $false, $false, $true, 1, $false, "Média Final = (N + Prova Recuperação)/2^l", 2)
if (-not $ok) { Write-Output "PHASE2 FAILED at index 7 (para=16 run=1)" }
$ok = $d.Content.Find.Execute("@@BLOCK8@@", $true, $false, $false, $false, $false, $true, 1, $false, "1- FOGLER, H.S. Elementos de engenharia das reações químicas. 3.ed. Rio de Janeiro: LTC Editora, 2009.^l^l2- LEVENSPIEL, O. Engenharia Das Reações Químicas, E ed (Blucher, São Paulo, 2000)^l3- VAN SANTEN, R.A.; Niemantsverdriet, J.W. Chemical kinetics and catalysis. New York: Plenum Press, 1995.^l4- Missen, R.W.; Mims, C.A.; Saville, B.A. Introduction to chemical reaction engineering and kinetics. New York: J. Wiley, 1999.^l5- Rothenberg, G. Catalysis: concepts and green applications. Weinheim: Wiley-VCH, 2008 Chichester.^l6- DENISOV, E.T.; Sarkisov, O.M.; Likhtenshtein, G.I. Chemical kinetics: fundamentals and new developments. Amsterdam: Elsevier, 2003.^l7- Hagen, J. Industrial catalysis: a practical approach. Weinheim: Wiley-VCH, 2006.^l8- Salmi, T.O.; Mikkola, J.; Warna, J.P. Chemical reaction engineering and reactor technology. Boca Raton: CRC Press/Taylor & Francis, 2011.^l9- Mortimer, M.; Taylor, P.G. Chemical kinetics and mechanism. Cambridge: Royal Society of Chemistry, 2002.^l10- FROMENT, G.F.; BISCHOFF, K.B. Chemical reactor analysis and design. 2nd. Ed. New York: John Wiley & Sons, 1990.^l11- HILL, C.G. An Introduction to chemical engineering kinetics and reactor design. New York: John Wiley&Sons, 1977.^l12- SMITH, J.M. Chemical engineering kinetics. 3rd. ed New York: McGraw-Hill,1981.^l13- DENBIGH, K.; TURNER, R. Introduction to chemical Reaction Design. Cambridge: Cambridge University Press, 1970.^l14 - AGUIAR, L. G. Problemas de cinética e reatores químicos. Curitiba: Appris Editora, 2023.^l", 2)
if (-not $ok) { Write-Output "PHASE2 FAILED at index 8 (para=16 run=3)" }
$ok = $d.Content.Find.Execute("@@BLOCK9@@", $true, $false, $false, $false, $false, $true, 1, $false, "5963230 - Leandro Gonçalves de Aguiar", 2)
if (-not $ok) { Write-Output "PHASE2 FAILED at index 9 (para=16 run=5)" }
$ok = $d.Content.Find.Execute("@@BLOCK10@@", $true, $false, $false, $false, $false, $true, 1, $false, "6310316 - Liana Alvares Rodrigues", 2)
if (-not $ok) { Write-Output "PHASE2 FAILED at index 10 (para=18 run=0)" }

Write-Output ("Total phase1 matches: " + $matches)
